$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Plain text columns (Coin name / Link) - direct assignment is safe
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'

# Numeric-looking text columns (Price / Volume) - force text to avoid numeric coercion
Set-TextValue $ws.Range("D2") '79.044.77'
Set-TextValue $ws.Range("E2") '  +3.35%  '
Set-TextValue $ws.Range("D3") '3.191.31'
Set-TextValue $ws.Range("E3") '  +5.25%  '
Set-TextValue $ws.Range("E4") '  -0.05%  '
Set-TextValue $ws.Range("D5") '206.27'
Set-TextValue $ws.Range("E5") '  +3.00%  '
Set-TextValue $ws.Range("D6") '632.76'
Set-TextValue $ws.Range("E6") '  +0.62%  '
Set-TextValue $ws.Range("E7") '  +0.02%  '
Set-TextValue $ws.Range("D8") '0.227'
Set-TextValue $ws.Range("E8") '  +11.24%  '
Set-TextValue $ws.Range("D9") '0.581'
Set-TextValue $ws.Range("E9") '  +5.38%  '
Set-TextValue $ws.Range("D10") '3.188.96'
Set-TextValue $ws.Range("E10") '  +5.17%  '
Set-TextValue $ws.Range("D11") '0.575'
Set-TextValue $ws.Range("E11") '  +31.43%  '
Set-TextValue $ws.Range("E12") '  +2.63%  '
Set-TextValue $ws.Range("D13") '5.43'
Set-TextValue $ws.Range("E13") '  +7.25%  '
Set-TextValue $ws.Range("D14") '3.769.05'
Set-TextValue $ws.Range("E14") '  +4.91%  '
Set-TextValue $ws.Range("D15") '0.0000224'
Set-TextValue $ws.Range("E15") '  +18.51%  '
Set-TextValue $ws.Range("D16") '31.50'
Set-TextValue $ws.Range("E16") '  +7.52%  '
Set-TextValue $ws.Range("D17") '78.919.54'
Set-TextValue $ws.Range("E17") '  +3.31%  '
Set-TextValue $ws.Range("D18") '3.182.43'
Set-TextValue $ws.Range("E18") '  +5.19%  '
Set-TextValue $ws.Range("D19") '14.42'
Set-TextValue $ws.Range("E19") '  +7.25%  '
Set-TextValue $ws.Range("D20") '9.36'
Set-TextValue $ws.Range("E20") '  +3.37%  '
Set-TextValue $ws.Range("D21") '427.53'
Set-TextValue $ws.Range("E21") '  +14.19%  '
Set-TextValue $ws.Range("D22") '2.84'
Set-TextValue $ws.Range("E22") '  +25.35%  '
Set-TextValue $ws.Range("E23") '  +12.90%  '
Set-TextValue $ws.Range("D24") '6.82'
Set-TextValue $ws.Range("E24") '  +5.59%  '
Set-TextValue $ws.Range("D25") '3.344.60'
Set-TextValue $ws.Range("E25") '  +4.31%  '
Set-TextValue $ws.Range("E26") '  +7.81%  '
Set-TextValue $ws.Range("D27") '75.71'
Set-TextValue $ws.Range("E27") '  +3.52%  '
Set-TextValue $ws.Range("D28") '10.89'
Set-TextValue $ws.Range("E28") '  +10.79%  '
Set-TextValue $ws.Range("D29") '1.01'
Set-TextValue $ws.Range("E29") '  +0.64%  '
Set-TextValue $ws.Range("D30") '0.0000114'
Set-TextValue $ws.Range("E30") '  +5.61%  '
Set-TextValue $ws.Range("D31") '0.995'
Set-TextValue $ws.Range("E31") '  -0.01%  '
Set-TextValue $ws.Range("D32") '8.79'
Set-TextValue $ws.Range("E32") '  +5.65%  '
Set-TextValue $ws.Range("D33") '1.48'
Set-TextValue $ws.Range("E33") '  +4.66%  '
Set-TextValue $ws.Range("D34") '513.81'
Set-TextValue $ws.Range("E34") '  +0.42%  '
Set-TextValue $ws.Range("D35") '1.99'
Set-TextValue $ws.Range("E35") '  +2.13%  '
Set-TextValue $ws.Range("D36") '22.84'
Set-TextValue $ws.Range("E36") '  +10.05%  '
Set-TextValue $ws.Range("D37") '0.125'
Set-TextValue $ws.Range("E37") '  +19.13%  '
Set-TextValue $ws.Range("D38") '0.134'
Set-TextValue $ws.Range("E38") '  +19.72%  '
Set-TextValue $ws.Range("D39") '1.00'
Set-TextValue $ws.Range("E39") '  +0.03%  '
Set-TextValue $ws.Range("D40") '0.398'
Set-TextValue $ws.Range("E40") '  +3.97%  '
Set-TextValue $ws.Range("D41") '164.11'
Set-TextValue $ws.Range("E41") '  +0.07%  '
Set-TextValue $ws.Range("D42") '19.97'
Set-TextValue $ws.Range("E42") '  -0.20%  '
Set-TextValue $ws.Range("D43") '192.14'
Set-TextValue $ws.Range("E43") '  -0.98%  '
Set-TextValue $ws.Range("D44") '1.00'
Set-TextValue $ws.Range("E44") '  -0.06%  '
Set-TextValue $ws.Range("D45") '5.41'
Set-TextValue $ws.Range("E45") '  +7.30%  '
Set-TextValue $ws.Range("D46") '0.802'
Set-TextValue $ws.Range("E46") '  +12.84%  '
Set-TextValue $ws.Range("E47") '  +7.61%  '
Set-TextValue $ws.Range("E48") '  +3.01%  '
Set-TextValue $ws.Range("D49") '42.43'
Set-TextValue $ws.Range("E49") '  -0.54%  '
Set-TextValue $ws.Range("D50") '2.49'
Set-TextValue $ws.Range("E50") '  +5.89%  '
Set-TextValue $ws.Range("D51") '0.621'
Set-TextValue $ws.Range("E51") '  +3.20%  '
